# 精制食用植物油.xlsx update
#
# 1. For each year, the "B" (Q2/半年) and "C" (Q3/三季度) period rows had
#    their data rows swapped (the row ordering in the sheet changes from
#    A,B,C,D to A,C,B,D per year-block; row labels travel with their data).
# 2. Columns F (产销率) and G (销售量) are dropped entirely — they were
#    redundant derived copies of columns B and E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Row($r1, $r2) {
    # Swap the full A:E contents of two rows (label + 4 data columns),
    # reading with Value2 (actual scalar) and writing back with Value.
    for ($col = 1; $col -le 5; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# (rowB, rowC) pairs — one per year present in the sheet.
Swap-Row 2 3    # 2004年B / 2004年C
Swap-Row 6 7    # 2005年B / 2005年C
Swap-Row 10 11  # 2006年B / 2006年C
Swap-Row 14 15  # 2007年B / 2007年C
Swap-Row 18 19  # 2008年B / 2008年C
Swap-Row 22 23  # 2009年B / 2009年C
Swap-Row 26 27  # 2010年B / 2010年C
Swap-Row 30 31  # 2011年B / 2011年C
Swap-Row 34 35  # 2012年B / 2012年C
Swap-Row 38 39  # 2016年B / 2016年C
Swap-Row 42 43  # 2017年B / 2017年C
Swap-Row 46 47  # 2018年B / 2018年C
Swap-Row 50 51  # 2019年B / 2019年C

# Drop the now-redundant 产销率 / 销售量 columns (F:G), shifting the used
# range back down to A1:E52.
$ws.Range("F1:G52").Delete()
